$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Return_with_prediction (G), return_pct_change (H) and, for row 2,
# mean_return_pct_change (I) with the recomputed values from the latest run.
$updates = @(
    @(2, 7, 0.03669687530071119),
    @(2, 8, -13.20236706568762),
    @(2, 9, 39.94334393601181),
    @(3, 7, 0.05888460316642855),
    @(3, 8, 20.91207584020855),
    @(4, 7, 0.01505732121492321),
    @(4, 8, 690.8690779321886),
    @(5, 7, 0.02503877817020134),
    @(5, 8, 559.5485685962095),
    @(6, 7, 0.03113033541913608),
    @(6, 8, -10.20850005435956),
    @(7, 7, 0.04926306059380776),
    @(7, 8, -7.381018269367439),
    @(8, 7, -0.03083676144967969),
    @(8, 8, -63.84427001137097),
    @(9, 7, -0.04391316466878215),
    @(9, 8, -103.9825270495194),
    @(10, 7, -0.05447779976807383),
    @(10, 8, 25.06632232947468),
    @(11, 7, -0.06047261554673258),
    @(11, 8, 34.27678402386258),
    @(12, 7, -0.2533697196301742),
    @(12, 8, -3.646924124319147),
    @(13, 7, -0.3698567553278699),
    @(13, 8, -34.58477510024577),
    @(14, 7, -0.06304516814824457),
    @(14, 8, -69.94592317147837),
    @(15, 7, -0.01123142833138899),
    @(15, 8, 67.69919682289331),
    @(16, 7, 0.1213001521279186),
    @(16, 8, -3.198542252622812),
    @(17, 7, 0.1249742704972536),
    @(17, 8, -10.89209781441065),
    @(18, 7, 0.1252435288847864),
    @(18, 8, 0.4095195414126414),
    @(19, 7, 0.1461661340715624),
    @(19, 8, 9.722931876804711),
    @(20, 7, 0.02861118657950151),
    @(20, 8, -16.67319239912537),
    @(21, 7, 0.05332588252388147),
    @(21, 8, -8.121969759385905),
    @(22, 7, -0.0794939282982165),
    @(22, 8, 0.4434123399933709),
    @(23, 7, -0.07851642164217705),
    @(23, 8, -25.53300979778504),
    @(24, 7, 0.09786712133948434),
    @(24, 8, -17.14216880862371),
    @(25, 7, 0.1395847881592913),
    @(25, 8, 10.6319668348486),
    @(26, 7, 0.05842503551548898),
    @(26, 8, 17.54568712224735),
    @(27, 7, 0.08934511227191817),
    @(27, 8, 3.078735070980489),
    @(28, 7, -0.0658456612588101),
    @(28, 8, -3.540622512446592),
    @(29, 7, -0.08097010043221027),
    @(29, 8, -13.7664374563126),
    @(30, 7, 0.08537414057522948),
    @(30, 8, 34.00824176033719),
    @(31, 7, 0.06993200905241485),
    @(31, 8, 15.43657409099624),
    @(32, 7, 0.09893828124753824),
    @(32, 8, 0.6862096984813666),
    @(33, 7, 0.1247381449753012),
    @(33, 8, 51.59335477067448),
    @(34, 7, -0.002150704869718027),
    @(34, 8, -108.2544119863903),
    @(35, 7, 0.008382497022382853),
    @(35, 8, 174.7762319353111),
    @(36, 7, 0.004428976118493459),
    @(36, 8, 722.9130146332789),
    @(37, 7, -0.006761602361657539),
    @(37, 8, 46.14099357258172),
    @(38, 7, 0.1041707524207047),
    @(38, 8, -2.877674478523505),
    @(39, 7, 0.1045586331169056),
    @(39, 8, 22.05841406790322),
    @(40, 7, 0.003767276113305323),
    @(40, 8, 26.83332615564793),
    @(41, 7, 0.01941374204941103),
    @(41, 8, 29.45302593698177),
    @(42, 7, 0.1013523672387794),
    @(42, 8, 0.4133761272055603),
    @(43, 7, 0.1222537702080459),
    @(43, 8, 1.755608399429374),
    @(44, 7, 0.03076492040577667),
    @(44, 8, -13.7937783762388),
    @(45, 7, 0.03572007648190318),
    @(45, 8, 118.2021774801854),
    @(46, 7, 0.05066437159375096),
    @(46, 8, 39.80806277835936),
    @(47, 7, 0.04562954200452892),
    @(47, 8, -9.537553383600367),
    @(48, 7, 0.04737124670390871),
    @(48, 8, 10.72885639368554),
    @(49, 7, 0.06568965953966789),
    @(49, 8, -5.449990852584079),
    @(50, 7, 0.01576581799100422),
    @(50, 8, -8.72410179436655),
    @(51, 7, 0.01599341387461513),
    @(51, 8, -17.85489863477662),
    @(52, 7, -0.1049607185442575),
    @(52, 8, -1.391283612921226),
    @(53, 7, -0.07165501097018118),
    @(53, 8, 22.4135081271865),
    @(54, 7, 0.06398983354280589),
    @(54, 8, -12.48973623977442),
    @(55, 7, 0.08671646937632152),
    @(55, 8, 39.97441409628797),
    @(56, 7, 0.02676707829122501),
    @(56, 8, -23.499060095036),
    @(57, 7, 0.03638097383518668),
    @(57, 8, 530.1353101877747),
    @(58, 7, 0.03974957792716755),
    @(58, 8, 58.93088083569654),
    @(59, 7, 0.03830675393947092),
    @(59, 8, 61.77777500529099),
    @(60, 7, 0.01000757853037881),
    @(60, 8, -69.15302522023239),
    @(61, 7, 0.04447619941784344),
    @(61, 8, 251.3684896269531),
    @(62, 7, 0.02872110653856923),
    @(62, 8, -52.41905084159801),
    @(63, 7, 0.04540136012176382),
    @(63, 8, 39.31287498552442),
    @(64, 7, 0.02677938367082512),
    @(64, 8, -33.92078157039081),
    @(65, 7, 0.06809490578783547),
    @(65, 8, 21.46373936593511),
    @(66, 7, 0.1058596686113273),
    @(66, 8, 13.15341969903649),
    @(67, 7, 0.125991945689305),
    @(67, 8, 9.134339259320463),
    @(68, 7, -0.05187881086540459),
    @(68, 8, -48.8610963018268),
    @(69, 7, -0.004755217170772516),
    @(69, 8, 77.59282754072117),
    @(70, 7, 0.07486329236143296),
    @(70, 8, -19.1857551111673),
    @(71, 7, 0.08014206986092634),
    @(71, 8, -12.13339818756186),
    @(72, 7, -0.03915423655426031),
    @(72, 8, 30.18273871152995),
    @(73, 7, -0.09160354971818271),
    @(73, 8, -24.18700302699491),
    @(74, 7, 0.111600460828608),
    @(74, 8, 11.65985628475841),
    @(75, 7, 0.1134939491778463),
    @(75, 8, 16.51874947637912),
    @(76, 7, -0.007559277526442381),
    @(76, 8, -129.5623802741607),
    @(77, 7, 0.0283558185390719),
    @(77, 8, 100.9762410567117),
    @(78, 7, 0.08798211001812058),
    @(78, 8, 36.87968165145399),
    @(79, 7, 0.06856901030124732),
    @(79, 8, -10.61722365494763),
    @(80, 7, -0.1502696939999564),
    @(80, 8, 9.262035087739259),
    @(81, 7, -0.1519912727742684),
    @(81, 8, 27.65490248316695),
    @(82, 7, 0.1491876000682133),
    @(82, 8, 30.07295296541117),
    @(83, 7, 0.1755919542025658),
    @(83, 8, -1.34297594715324),
    @(84, 7, 0.08816241887359896),
    @(84, 8, 269.8452758920145),
    @(85, 7, 0.04485108812233483),
    @(85, 8, -27.16134260877758),
)

foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}